# Rebuild the "baptista.doi" DOI list: merge in newly-added DOIs and
# re-sort the whole list A-Z, turning every DOI cell into a hyperlink
# (https://doi.org/<doi>) styled blue + underlined, same as Google
# Sheets' "Insert link" treatment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dois = @(
      "10.1007/s00170-023-12686-2",
      "10.1007/s11665-022-06795-4",
      "10.1007/s40430-014-0175-6",
      "10.1007/s40430-019-1821-9",
      "10.1007/s40430-023-04540-z",
      "10.1007/s40430-024-05069-5",
      "10.1007/s42247-024-00746-6",
      "10.1016/j.ceramint.2012.08.093",
      "10.1016/j.engfailanal.2006.11.053",
      "10.1016/j.ijfatigue.2004.01.011",
      "10.1016/j.ijfatigue.2016.04.018",
      "10.1016/j.jmbbm.2014.06.011",
      "10.1016/j.jmrt.2013.10.005",
      "10.1016/j.jmrt.2021.01.038",
      "10.1016/j.matdes.2008.08.010",
      "10.1016/j.mechmat.2012.04.003",
      "10.1016/j.msea.2006.10.124",
      "10.1016/j.msea.2025.148040",
      "10.1016/j.proeng.2010.03.178",
      "10.1016/j.proeng.2010.03.205",
      "10.1016/j.proeng.2010.03.206",
      "10.1016/j.proeng.2011.04.202",
      "10.1016/j.proeng.2011.04.220",
      "10.1016/j.proeng.2011.04.339",
      "10.1016/j.proeng.2011.04.342",
      "10.1016/j.proeng.2011.04.343",
      "10.1016/j.prostr.2019.08.043",
      "10.1016/j.prostr.2019.08.054",
      "10.1051/matecconf/202032111095",
      "10.1088/2053-1591/ab6c8f",
      "10.1177/0021998310376098",
      "10.1243/14644207JMDA122",
      "10.1590/0366-69132022683863210",
      "10.1590/0370-44672018720157",
      "10.1590/0370-44672019740154",
      "10.1590/1980-5373-MR-2022-0315",
      "10.1590/s0366-69132014000100006",
      "10.17563/rbav.v34i1.939",
      "10.17563/rbav.v34i2.956",
      "10.20906/CPS/COB-2015-0858",
      "10.20906/CPS/COB-2015-0934",
      "10.20906/CPS/COB-2015-1461",
      "10.26678/ABCM.COBEM2023.COB2023-0535",
      "10.3390/met11122034",
      "10.3390/met9050589",
      "10.4028/www.scientific.net/AMR.891-892.1335",
      "10.4028/www.scientific.net/AMR.891-892.1507",
      "10.4028/www.scientific.net/AMR.891-892.1767",
      "10.4028/www.scientific.net/AMR.891-892.1785",
      "10.4028/www.scientific.net/AMR.891-892.615",
      "10.4028/www.scientific.net/KEM.713.30",
      "10.4028/www.scientific.net/MSF.636-637.47",
      "10.4028/www.scientific.net/MSF.727-728.945",
      "10.4028/www.scientific.net/amr.1135.167",
      "10.4028/www.scientific.net/kem.754.248",
      "10.4028/www.scientific.net/kem.754.43",
      "10.4028/www.scientific.net/msf.805.204",
      "10.4028/www.scientific.net/msf.869.503",
      "10.5028/jatm.v5i2.219"
)

$firstRow = 2
$lastRow = $firstRow + $dois.Length - 1

# Write every DOI in one shot as a single-column array so the sheet ends
# up with exactly one header row + 59 sorted data rows.
$target = $ws.Range($ws.Cells.Item($firstRow, 1), $ws.Cells.Item($lastRow, 1))
$arr = New-Object 'object[,]' $dois.Length,1
for ($i = 0; $i -lt $dois.Length; $i++) {
    $arr[$i, 0] = $dois[$i]
}
$target.Value2 = $arr

# Blue + underlined hyperlink look for the whole range (cell-level font).
$target.Font.Underline = $true
$target.Font.Color = 255 * 65536

# Per-cell hyperlink + matching rich-text run colour inside the cell text
# (Google Sheets link colour #1155CC) so the shared-string run carries its
# own <rPr> too, not just the cell style.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $doi = $cell.Value2
    $ws.Hyperlinks.Add($cell, ("https://doi.org/" + $doi))

    $len = $doi.Length
    if ($len -gt 1) {
        $cell.Characters(1, $len - 1).Font.Color = 13391121
        $cell.Characters(1, $len - 1).Font.Underline = $true
    }
    $cell.Characters($len, 1).Font.Color = 13391121
    $cell.Characters($len, 1).Font.Underline = $true
}
